$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set text format on the affected Price/Coin/Link/Volume columns so that
# numeric-looking strings (e.g. "1.00", "4.90") are kept as literal text
# instead of being coerced to numbers and losing trailing zeros.
$ws.Range("B2:E51").NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '46.907.00'
$ws.Range('E2').Value = '  +5.30%  '

# Row 3
$ws.Range('D3').Value = '2.329.10'
$ws.Range('E3').Value = '  +3.67%  '

# Row 4
$ws.Range('E4').Value = '  -0.78%  '

# Row 5
$ws.Range('D5').Value = '306.19'
$ws.Range('E5').Value = '  +0.25%  '

# Row 6
$ws.Range('D6').Value = '97.89'
$ws.Range('E6').Value = '  +3.56%  '

# Row 7
$ws.Range('D7').Value = '0.578'
$ws.Range('E7').Value = '  +1.38%  '

# Row 8
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.66%  '

# Row 9
$ws.Range('D9').Value = '0.537'
$ws.Range('E9').Value = '  +3.69%  '

# Row 10
$ws.Range('D10').Value = '35.88'
$ws.Range('E10').Value = '  +3.21%  '

# Row 11
$ws.Range('D11').Value = '0.0807'
$ws.Range('E11').Value = '  +0.71%  '

# Row 12
$ws.Range('D12').Value = '7.45'
$ws.Range('E12').Value = '  +3.38%  '

# Row 14
$ws.Range('D14').Value = '2.690.26'
$ws.Range('E14').Value = '  +3.84%  '

# Row 15
$ws.Range('D15').Value = '2.339.12'
$ws.Range('E15').Value = '  +4.42%  '

# Row 16
$ws.Range('D16').Value = '14.10'
$ws.Range('E16').Value = '  +4.10%  '

# Row 17
$ws.Range('D17').Value = '0.829'
$ws.Range('E17').Value = '  -0.37%  '

# Row 18
$ws.Range('D18').Value = '46.683.72'
$ws.Range('E18').Value = '  +5.25%  '

# Row 19
$ws.Range('D19').Value = '13.14'
$ws.Range('E19').Value = '  +11.87%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0945'
$ws.Range('E20').Value = '  +0.74%  '

# Row 21
$ws.Range('D21').Value = '6.18'
$ws.Range('E21').Value = '  -0.07%  '

# Row 22
$ws.Range('D22').Value = '66.66'
$ws.Range('E22').Value = '  +2.20%  '

# Row 23
$ws.Range('D23').Value = '246.44'
$ws.Range('E23').Value = '  +3.90%  '

# Row 24
$ws.Range('D24').Value = '2.97'
$ws.Range('E24').Value = '  +0.87%  '

# Row 25
$ws.Range('E25').Value = '  -0.07%  '

# Row 26
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.26%  '

# Row 27
$ws.Range('D27').Value = '41.93'
$ws.Range('E27').Value = '  +12.86%  '

# Row 28
$ws.Range('D28').Value = '2.29'
$ws.Range('E28').Value = '  -2.53%  '

# Row 29
$ws.Range('E29').Value = '  +0.85%  '

# Row 30
$ws.Range('D30').Value = '20.12'
$ws.Range('E30').Value = '  +0.82%  '

# Row 31
$ws.Range('D31').Value = '5.70'
$ws.Range('E31').Value = '  -3.52%  '

# Row 32
$ws.Range('D32').Value = '151.55'
$ws.Range('E32').Value = '  +1.32%  '

# Row 33
$ws.Range('E33').Value = '  +3.01%  '

# Row 34
$ws.Range('E34').Value = '  +0.75%  '

# Row 35
$ws.Range('D35').Value = '3.02'
$ws.Range('E35').Value = '  -5.45%  '

# Row 36
$ws.Range('E36').Value = '  +0.47%  '

# Row 37
$ws.Range('E37').Value = '  -2.04%  '

# Row 38
$ws.Range('D38').Value = '1.80'
$ws.Range('E38').Value = '  -3.50%  '

# Row 39
$ws.Range('E39').Value = '  +6.43%  '

# Row 40
$ws.Range('E40').Value = '  +6.67%  '

# Row 41
$ws.Range('E41').Value = '  +1.41%  '

# Row 42
$ws.Range('D42').Value = '13.86'
$ws.Range('E42').Value = '  -8.88%  '

# Row 43
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  -0.74%  '

# Row 44
$ws.Range('D44').Value = '1.96'
$ws.Range('E44').Value = '  +10.39%  '

# Row 45
$ws.Range('D45').Value = '1.792.85'
$ws.Range('E45').Value = '  -0.61%  '

# Row 46
$ws.Range('E46').Value = '  +6.11%  '

# Row 47
$ws.Range('B47').Value = 'BitcoinSV'
$ws.Range('C47').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D47').Value = '80.91'
$ws.Range('E47').Value = '  -1.20%  '

# Row 48
$ws.Range('B48').Value = 'ordi'
$ws.Range('C48').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D48').Value = '73.57'
$ws.Range('E48').Value = '  +7.11%  '

# Row 49
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').Value = '4.90'
$ws.Range('E49').Value = '  +1.47%  '

# Row 50
$ws.Range('D50').Value = '98.00'
$ws.Range('E50').Value = '  -0.45%  '

# Row 51
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').Value = '55.28'
$ws.Range('E51').Value = '  +2.62%  '
